$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ESTADOS UNIDOS"
$ws.Range("B2").Value = "NA"
$ws.Range("C2").Value = "'94500"

$ws.Range("A3").Value = "ESTADOS UNIDOS"
$ws.Range("B3").Value = "EU"
$ws.Range("C3").Value = "504500"

$ws.Range("A4").Value = "RUSIA"
$ws.Range("B4").Value = "EUOR"
$ws.Range("C4").Value = "760000"

$ws.Range("A5").Value = "INGLATERRA"
$ws.Range("B5").Value = "EUOC"
$ws.Range("C5").Value = "709000"
